# Auto-generated edit script to update cryptos.xlsx price/volume/coin data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '37.246.41'),
    @('E2', '  +2.45%  '),
    @('D3', '2.061.33'),
    @('E3', '  +2.02%  '),
    @('E4', '  -0.08%  '),
    @('D5', '232.14'),
    @('E5', '  +0.51%  '),
    @('D6', '0.618'),
    @('E6', '  +3.33%  '),
    @('E7', '  -0.10%  '),
    @('D8', '57.70'),
    @('E8', '  +6.08%  '),
    @('E9', '  +3.52%  '),
    @('B10', 'Dogecoin'),
    @('C10', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'),
    @('D10', '0.0757'),
    @('E10', '  +1.77%  '),
    @('B11', 'TRON'),
    @('C11', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'),
    @('D11', '0.101'),
    @('E11', '  +1.22%  '),
    @('B12', 'WrappedliquidstakedEther2.0'),
    @('C12', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'),
    @('D12', '2.365.32'),
    @('E12', '  +1.88%  '),
    @('B13', 'Chainlink'),
    @('C13', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'),
    @('D13', '14.33'),
    @('E13', '  +1.03%  '),
    @('B14', 'Avalanche'),
    @('C14', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'),
    @('D14', '20.80'),
    @('E14', '  +4.15%  '),
    @('B15', 'Polygon'),
    @('C15', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'),
    @('D15', '0.772'),
    @('E15', '  +2.18%  '),
    @('B16', 'Polkadot'),
    @('C16', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'),
    @('D16', '5.13'),
    @('E16', '  +1.00%  '),
    @('B17', 'WrappedEther'),
    @('C17', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'),
    @('D17', '2.061.10'),
    @('E17', '  +2.14%  '),
    @('B18', 'WrappedBTC'),
    @('C18', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'),
    @('D18', '37.188.28'),
    @('E18', '  +1.71%  '),
    @('B19', 'Uniswap'),
    @('C19', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'),
    @('D19', '6.37'),
    @('E19', '  +16.71%  '),
    @('B20', 'Litecoin'),
    @('C20', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'),
    @('D20', '69.09'),
    @('E20', '  +2.53%  '),
    @('B21', 'ShibaInu'),
    @('C21', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'),
    @('D21', '0.0₃0809'),
    @('E21', '  +2.08%  '),
    @('B22', 'BitcoinCash'),
    @('C22', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'),
    @('D22', '224.58'),
    @('E22', '  +2.11%  '),
    @('B23', 'Dai'),
    @('C23', 'https://coinranking.com/coin/MoTuySvg7+dai-dai'),
    @('D23', '1.00'),
    @('E23', '  -0.15%  '),
    @('B24', 'Toncoin'),
    @('C24', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'),
    @('D24', '2.43'),
    @('E24', '  +2.44%  '),
    @('B25', 'PancakeSwap'),
    @('C25', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'),
    @('D25', '2.37'),
    @('E25', '  +0.63%  '),
    @('B26', 'Monero'),
    @('C26', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'),
    @('D26', '165.61'),
    @('E26', '  +1.88%  '),
    @('B27', 'ImmutableX'),
    @('C27', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'),
    @('D27', '1.45'),
    @('E27', '  +7.77%  '),
    @('B28', 'Cosmos'),
    @('C28', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'),
    @('D28', '8.77'),
    @('E28', '  +2.09%  '),
    @('B29', 'EthereumClassic'),
    @('C29', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'),
    @('D29', '19.05'),
    @('E29', '  +1.67%  '),
    @('D30', '0.126'),
    @('E30', '  -1.17%  '),
    @('B31', 'Stellar'),
    @('C31', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'),
    @('D31', '0.117'),
    @('E31', '  +0.21%  '),
    @('B32', 'Filecoin'),
    @('C32', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'),
    @('D32', '4.45'),
    @('E32', '  +2.32%  '),
    @('B33', 'LidoDAOToken'),
    @('C33', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'),
    @('D33', '2.54'),
    @('E33', '  +3.38%  '),
    @('D34', '0.0612'),
    @('E34', '  +2.38%  '),
    @('B35', 'InternetComputer(DFINITY)'),
    @('C35', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'),
    @('D35', '4.52'),
    @('E35', '  +6.81%  '),
    @('B36', 'BinanceUSD'),
    @('C36', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'),
    @('D36', '1.00'),
    @('E36', '  -0.18%  '),
    @('B37', 'RenderToken'),
    @('C37', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @('D37', '3.27'),
    @('E37', '  +0.18%  '),
    @('E38', '  -1.04%  '),
    @('B39', 'THORChain'),
    @('C39', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'),
    @('D39', '5.78'),
    @('E39', '  +0.13%  '),
    @('B40', 'FTXToken'),
    @('C40', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'),
    @('D40', '4.62'),
    @('E40', '  +17.12%  '),
    @('B41', 'HuobiToken'),
    @('C41', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'),
    @('D41', '2.95'),
    @('E41', '  +1.27%  '),
    @('B42', 'Maker'),
    @('C42', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'),
    @('D42', '1.487.76'),
    @('E42', '  +1.86%  '),
    @('B43', 'Aave'),
    @('C43', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'),
    @('D43', '96.29'),
    @('E43', '  +3.88%  '),
    @('B44', 'Cronos'),
    @('C44', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'),
    @('D44', '0.0930'),
    @('E44', '  +0.18%  '),
    @('E45', '  +5.23%  '),
    @('B46', 'VeChain'),
    @('C46', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('D46', '0.0210'),
    @('E46', '  +3.90%  '),
    @('B47', 'InjectiveProtocol'),
    @('C47', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'),
    @('D47', '15.29'),
    @('E47', '  -0.92%  '),
    @('B48', 'ARBITRUM'),
    @('C48', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'),
    @('D48', '1.02'),
    @('E48', '  +2.69%  '),
    @('B49', 'FraxShare'),
    @('C49', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'),
    @('D49', '7.14'),
    @('E49', '  +4.55%  '),
    @('E50', '  +2.57%  '),
    @('B51', 'RocketPoolETH'),
    @('C51', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'),
    @('D51', '2.252.24'),
    @('E51', '  +1.84%  ')
)

foreach ($u in $updates) {
    $cell = $ws.Range($u[0])
    $cell.NumberFormat = "@"
    $cell.Value = $u[1]
    $cell.ClearFormats()
}

